# Updates the cryptos price/volume table with the latest scraped values
# (GitHub Actions crypto tracker refresh). For numeric-looking "Price"
# strings we force text storage (NumberFormat "@") so Excel doesn't
# reinterpret/round them as numbers, then restore the default "Normal"
# style so no stray cell formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.300.59'
$ws.Range('E2').Value = '  +0.43%  '
$ws.Range('D3').Value = '1.589.63'
$ws.Range('E4').Value = '  -0.24%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '212.04'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  +1.12%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '0.502'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  +0.83%  '
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('E9').Value = '  -0.16%  '
$ws.Range('E10').Value = '  -0.70%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '0.0846'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  -0.12%  '
$ws.Range('D12').Value = '1.813.59'
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').Value = '1.589.92'
$ws.Range('E13').Value = '  +0.35%  '
$ws.Range('E14').Value = '  +0.31%  '
$ws.Range('E15').Value = '  +0.96%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '64.44'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').Value = '26.312.16'
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('D18').Value = '0.0₃0732'
$ws.Range('E18').Value = '  -0.18%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '7.50'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  +3.46%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '212.64'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  +3.06%  '
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('E22').Value = '  +0.83%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '8.99'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  +1.54%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '2.14'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  -2.45%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '145.08'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  +0.16%  '
$ws.Range('E26').Value = '  -0.16%  '
$ws.Range('E27').Value = '  +0.27%  '
$ws.Range('E28').Value = '  -0.26%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '15.19'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('E30').Value = '  -0.71%  '
$ws.Range('E31').Value = '  +0.75%  '
$ws.Range('E32').Value = '  +0.01%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '2.98'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  +1.24%  '
$ws.Range('D34').Value = '1.341.67'
$ws.Range('E34').Value = '  +4.66%  '
$ws.Range('E35').Value = '  -0.84%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '0.600'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  -0.63%  '
$ws.Range('E37').Value = '  +0.24%  '
$ws.Range('E38').Value = '  -0.10%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '1.04'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  -15.53%  '
$ws.Range('E40').Value = '  +0.55%  '
$ws.Range('E41').Value = '  +3.89%  '
$ws.Range('E42').Value = '  -0.19%  '
$ws.Range('E43').Value = '  +0.23%  '
$ws.Range('E44').Value = '  -0.75%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.725.13'
$ws.Range('E45').Value = '  +0.35%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '61.76'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  -0.73%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '87.99'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  -0.84%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0103'
$ws.Range('E48').Value = '  +6.95%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '1.49'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  -4.04%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '0.0981'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -2.60%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '0.0503'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -0.74%  '
